$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: updated odds values
$ws.Range("G2").Value = 2.18
$ws.Range("H2").Value = 2.85
$ws.Range("I2").Value = 3.7
$ws.Range("J2").Value = 2.85
$ws.Range("K2").Value = 1.88
$ws.Range("L2").Value = 4.45
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 5.1
$ws.Range("O2").Value = 1.55
$ws.Range("P2").Value = 2.3
$ws.Range("Q2").Value = 2.62
$ws.Range("R2").Value = 1.44
$ws.Range("S2").Value = 4.75
$ws.Range("T2").Value = 1.15
$ws.Range("U2").Value = 1.57
$ws.Range("V2").Value = 2.25
$ws.Range("W2").Value = 2.18
$ws.Range("X2").Value = 1.62
$ws.Range("Y2").Value = 5.5
$ws.Range("Z2").Value = 9
$ws.Range("AA2").Value = 9.5
$ws.Range("AB2").Value = 21
$ws.Range("AC2").Value = 22
$ws.Range("AD2").Value = 45
$ws.Range("AE2").Value = 5.1
$ws.Range("AF2").Value = 5.8
$ws.Range("AG2").Value = 19
$ws.Range("AH2").Value = 120
$ws.Range("AI2").Value = 7.5
$ws.Range("AJ2").Value = 17.5
$ws.Range("AK2").Value = 14
$ws.Range("AL2").Value = 60
$ws.Range("AM2").Value = 45
$ws.Range("AN2").Value = 70

# Row 3: updated odds values
$ws.Range("G3").Value = 1.9
$ws.Range("H3").Value = 3.8
$ws.Range("I3").Value = 3.35
$ws.Range("J3").Value = 2.42
$ws.Range("K3").Value = 2.3
$ws.Range("L3").Value = 3.7
$ws.Range("O3").Value = 1.19
$ws.Range("P3").Value = 3.7
$ws.Range("Q3").Value = 1.57
$ws.Range("R3").Value = 2.1
$ws.Range("S3").Value = 2.37
$ws.Range("T3").Value = 1.45
$ws.Range("W3").Value = 1.55
$ws.Range("X3").Value = 2.15
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 10.5
$ws.Range("AA3").Value = 8.5
$ws.Range("AB3").Value = 17
$ws.Range("AC3").Value = 13.5
$ws.Range("AD3").Value = 21
$ws.Range("AE3").Value = 15
$ws.Range("AF3").Value = 7.6
$ws.Range("AH3").Value = 45
$ws.Range("AI3").Value = 13
$ws.Range("AJ3").Value = 20
$ws.Range("AK3").Value = 11.75
$ws.Range("AL3").Value = 45
$ws.Range("AM3").Value = 26
$ws.Range("AO3").Value = 300

# Row 4: updated odds values
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 13
$ws.Range("Q4").Value = 1.63
$ws.Range("W4").Value = 1.67
$ws.Range("X4").Value = 2.1
$ws.Range("AF4").Value = 8

# Row 5: updated odds values
$ws.Range("G5").Value = 1.36
$ws.Range("H5").Value = 5.25
$ws.Range("I5").Value = 6.5
$ws.Range("J5").Value = 1.8
$ws.Range("L5").Value = 6
$ws.Range("O5").Value = 1.11
$ws.Range("P5").Value = 6
$ws.Range("Q5").Value = 1.37
$ws.Range("R5").Value = 2.75
$ws.Range("W5").Value = 1.62
$ws.Range("X5").Value = 2.2
$ws.Range("Y5").Value = 11
$ws.Range("Z5").Value = 9
$ws.Range("AB5").Value = 10
$ws.Range("AG5").Value = 17
$ws.Range("AL5").Value = 67

# Row 6: updated odds values
$ws.Range("Q6").Value = 1.41
$ws.Range("R6").Value = 2.62

# Row 7: updated odds values
$ws.Range("G7").Value = 1.88
$ws.Range("H7").Value = 3.2
$ws.Range("I7").Value = 4.2
$ws.Range("J7").Value = 2.4
$ws.Range("L7").Value = 4.55
$ws.Range("N7").Value = 6.4
$ws.Range("S7").Value = 3.55
$ws.Range("U7").Value = 1.39
$ws.Range("V7").Value = 2.75
$ws.Range("W7").Value = 1.9
$ws.Range("Y7").Value = 6.3
$ws.Range("Z7").Value = 8.25
$ws.Range("AA7").Value = 8.25
$ws.Range("AB7").Value = 16
$ws.Range("AC7").Value = 16
$ws.Range("AD7").Value = 30
$ws.Range("AE7").Value = 6.4
$ws.Range("AF7").Value = 6.2
$ws.Range("AI7").Value = 10.25
$ws.Range("AJ7").Value = 23
$ws.Range("AK7").Value = 14
$ws.Range("AL7").Value = 70
$ws.Range("AM7").Value = 45
$ws.Range("AN7").Value = 50
